$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 52700
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 102400
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 102400
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -102896

$ws.Range("H67").Value = 52700
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 102400
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 102400
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -104116

$ws.Range("H76").Value = 27780846
$ws.Range("I76").Value = 32261160
$ws.Range("J76").Value = 2900
$ws.Range("K76").Value = 32261160
$ws.Range("L76").Value = 2900
$ws.Range("M76").Value = -32260845
$ws.Range("N76").Value = -3530

$ws.Range("H79").Value = 27780846
$ws.Range("I79").Value = 32261160
$ws.Range("J79").Value = 2900
$ws.Range("K79").Value = 32261160
$ws.Range("L79").Value = 2900
$ws.Range("M79").Value = -32260068
$ws.Range("N79").Value = -5084

$ws.Range("H88").Value = 1000.25
$ws.Range("I88").Value = 901
$ws.Range("J88").Value = 1033.3334
$ws.Range("K88").Value = 901
$ws.Range("L88").Value = 1033.3334
$ws.Range("M88").Value = -495
$ws.Range("N88").Value = -1845.3334

$ws.Range("H91").Value = 1000.25
$ws.Range("I91").Value = 901
$ws.Range("J91").Value = 1033.3334
$ws.Range("K91").Value = 901
$ws.Range("L91").Value = 1033.3334
$ws.Range("M91").Value = 503
$ws.Range("N91").Value = -3841.3334

$ws.Range("H129").Value = 852.75
$ws.Range("I129").Value = 670.5714
$ws.Range("J129").Value = 994.44446
$ws.Range("K129").Value = 2011.7142
$ws.Range("L129").Value = 2983.33338
$ws.Range("M129").Value = 2988.2858
$ws.Range("N129").Value = -12983.33338

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21384.879
$ws.Range("I32").Value = 13490.295
$ws.Range("J32").Value = 34057.76
$ws.Range("K32").Value = 13490.295
$ws.Range("L32").Value = 34057.76
$ws.Range("M32").Value = -13203.295
$ws.Range("N32").Value = -34631.76

$ws.Range("H61").Value = 1497.8334
$ws.Range("I61").Value = 1291.8235
$ws.Range("K61").Value = 1291.8235
$ws.Range("M61").Value = -1079.8235

$ws.Range("H63").Value = 3200.5
$ws.Range("I63").Value = 3000.5557
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 3000.5557
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -2314.5557

$ws.Range("H66").Value = 3200.5
$ws.Range("I66").Value = 3000.5557
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 15002.7785
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -11570.7785

$ws.Range("H132").Value = 206255.58
$ws.Range("I132").Value = 8195.866
$ws.Range("J132").Value = 503345.16
$ws.Range("K132").Value = 24587.598
$ws.Range("L132").Value = 1510035.48
$ws.Range("M132").Value = -22057.598
$ws.Range("N132").Value = -1515095.48

$ws.Range("H136").Value = 1497.8334
$ws.Range("I136").Value = 1291.8235
$ws.Range("K136").Value = 3875.4705
$ws.Range("M136").Value = -1325.4705

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 23863.666
$ws.Range("J81").Value = 23863.666
$ws.Range("L81").Value = 23863.666
$ws.Range("N81").Value = -25985.666

$ws.Range("H84").Value = 23863.666
$ws.Range("J84").Value = 23863.666
$ws.Range("L84").Value = 71590.99800000001
$ws.Range("N84").Value = -82198.99800000001

$ws.Range("H105").Value = 1138317.2
$ws.Range("I105").Value = 1991105
$ws.Range("J105").Value = 1267
$ws.Range("K105").Value = 1991105
$ws.Range("L105").Value = 1267
$ws.Range("M105").Value = -1989358
$ws.Range("N105").Value = -4761

$ws.Range("H134").Value = 38502360
$ws.Range("I134").Value = 2441.8667
$ws.Range("J134").Value = 91002260
$ws.Range("K134").Value = 7325.6001
$ws.Range("L134").Value = 273006780
$ws.Range("M134").Value = -4790.6001
$ws.Range("N134").Value = -273011850

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19858.285
$ws.Range("I31").Value = 24190.113
$ws.Range("K31").Value = 24190.113
$ws.Range("M31").Value = -23895.113

$ws.Range("H34").Value = 19858.285
$ws.Range("I34").Value = 24190.113
$ws.Range("K34").Value = 24190.113
$ws.Range("M34").Value = -23988.113

$ws.Range("H62").Value = 3493.25
$ws.Range("I62").Value = 2502.5
$ws.Range("K62").Value = 2502.5
$ws.Range("M62").Value = -1878.5

$ws.Range("H65").Value = 3493.25
$ws.Range("I65").Value = 2502.5
$ws.Range("K65").Value = 12512.5
$ws.Range("M65").Value = -9392.5

$ws.Range("H141").Value = 42732.707
$ws.Range("J141").Value = 45803.867
$ws.Range("L141").Value = 45803.867
$ws.Range("N141").Value = -56163.867

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 124267.02
$ws.Range("I131").Value = 427.22223
$ws.Range("J131").Value = 159649.83
$ws.Range("K131").Value = 1281.66669
$ws.Range("L131").Value = 478949.49
$ws.Range("M131").Value = 3758.33331
$ws.Range("N131").Value = -489029.49

$ws.Range("H137").Value = 7319076
$ws.Range("I137").Value = 20834814
$ws.Range("J137").Value = 3714879.5
$ws.Range("K137").Value = 62504442
$ws.Range("L137").Value = 11144638.5
$ws.Range("M137").Value = -62499342
$ws.Range("N137").Value = -11154838.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4027.818
$ws.Range("I70").Value = 3799.8
$ws.Range("J70").Value = 4217.8335
$ws.Range("K70").Value = 3799.8
$ws.Range("L70").Value = 4217.8335
$ws.Range("M70").Value = -3529.8
$ws.Range("N70").Value = -4757.8335

$ws.Range("H73").Value = 4027.818
$ws.Range("I73").Value = 3799.8
$ws.Range("J73").Value = 4217.8335
$ws.Range("K73").Value = 3799.8
$ws.Range("L73").Value = 4217.8335
$ws.Range("M73").Value = -2863.8
$ws.Range("N73").Value = -6089.8335

$ws.Range("H132").Value = 35362.766
$ws.Range("I132").Value = 1781.15
$ws.Range("J132").Value = 102526
$ws.Range("K132").Value = 5343.450000000001
$ws.Range("L132").Value = 307578
$ws.Range("M132").Value = -2813.450000000001
$ws.Range("N132").Value = -312638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 28098.564
$ws.Range("I40").Value = 1775.4166
$ws.Range("K40").Value = 1775.4166
$ws.Range("M40").Value = -1639.4166

$ws.Range("H68").Value = 2176.4546
$ws.Range("I68").Value = 1356.5834
$ws.Range("J68").Value = 3160.3
$ws.Range("K68").Value = 1356.5834
$ws.Range("L68").Value = 3160.3
$ws.Range("M68").Value = -607.5834
$ws.Range("N68").Value = -4658.3

$ws.Range("H71").Value = 2176.4546
$ws.Range("I71").Value = 1356.5834
$ws.Range("J71").Value = 3160.3
$ws.Range("K71").Value = 6782.916999999999
$ws.Range("L71").Value = 15801.5
$ws.Range("M71").Value = -3038.916999999999
$ws.Range("N71").Value = -23289.5

$ws.Range("H122").Value = 2587.2942
$ws.Range("I122").Value = 2626
$ws.Range("J122").Value = 2494.4
$ws.Range("K122").Value = 7878
$ws.Range("L122").Value = 7483.200000000001
$ws.Range("M122").Value = -5428
$ws.Range("N122").Value = -12383.2

$ws.Range("H136").Value = 419000.16
$ws.Range("I136").Value = 835582
$ws.Range("J136").Value = 2418.3333
$ws.Range("K136").Value = 2506746
$ws.Range("L136").Value = 7254.999899999999
$ws.Range("M136").Value = -2504196
$ws.Range("N136").Value = -12354.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4314.6177
$ws.Range("J62").Value = 4333.6665
$ws.Range("L62").Value = 4333.6665
$ws.Range("N62").Value = -5581.6665

$ws.Range("H65").Value = 4314.6177
$ws.Range("J65").Value = 4333.6665
$ws.Range("L65").Value = 21668.3325
$ws.Range("N65").Value = -27908.3325

$ws.Range("H81").Value = 2383.5652
$ws.Range("I81").Value = 2077.1428
$ws.Range("J81").Value = 2860.2222
$ws.Range("K81").Value = 4154.2856
$ws.Range("L81").Value = 5720.4444
$ws.Range("M81").Value = -3093.2856
$ws.Range("N81").Value = -7842.4444

$ws.Range("H84").Value = 2383.5652
$ws.Range("I84").Value = 2077.1428
$ws.Range("J84").Value = 2860.2222
$ws.Range("K84").Value = 20771.428
$ws.Range("L84").Value = 28602.222
$ws.Range("M84").Value = -15467.428
$ws.Range("N84").Value = -39210.222

$ws.Range("H122").Value = 1487.3636
$ws.Range("I122").Value = 1026.0952
$ws.Range("J122").Value = 2294.5833
$ws.Range("K122").Value = 3078.2856
$ws.Range("L122").Value = 6883.749899999999
$ws.Range("M122").Value = -628.2856000000002
$ws.Range("N122").Value = -11783.7499

$ws.Range("H136").Value = 1698398.8
$ws.Range("I136").Value = 2382571
$ws.Range("J136").Value = 558111.4
$ws.Range("K136").Value = 7147713
$ws.Range("L136").Value = 1674334.2
$ws.Range("M136").Value = -7145163
$ws.Range("N136").Value = -1679434.2
